$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 (B5:AH5) values to 2 decimal places (custom accuracy).
$ws.Range("B5").Value  = 1.02
$ws.Range("C5").Value  = 0.46
$ws.Range("D5").Value  = 0.01
$ws.Range("E5").Value  = 0.34
$ws.Range("F5").Value  = 0.45
$ws.Range("H5").Value  = 11.07
$ws.Range("I5").Value  = 0.2
$ws.Range("J5").Value  = 0.5
$ws.Range("K5").Value  = 0.79
$ws.Range("O5").Value  = 0.14
$ws.Range("P5").Value  = 1.12
$ws.Range("Q5").Value  = 0.17
$ws.Range("R5").Value  = 0.18
$ws.Range("U5").Value  = 1.8
$ws.Range("V5").Value  = 0.51
$ws.Range("W5").Value  = 1.61
$ws.Range("X5").Value  = 0.46
$ws.Range("Z5").Value  = 4.8
$ws.Range("AA5").Value = 0.28
$ws.Range("AB5").Value = 0.09
$ws.Range("AC5").Value = 0.13
$ws.Range("AD5").Value = 0.76
$ws.Range("AE5").Value = 0.52
$ws.Range("AF5").Value = 10.62
$ws.Range("AG5").Value = 0.04
$ws.Range("AH5").Value = 0.3

# Remove the last data row (row 6) entirely - data trimmed to 1000 rows overall.
$ws.Rows("6:6").Delete()
